$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 74.41666666666666
$ws.Range("C2").Value = 74.82552934298489
$ws.Range("D2").Value = 74.41666666666666
$ws.Range("E2").Value = 74.3108953461858
$ws.Range("F2").Value = 0.7928402777777778
$ws.Range("G2").Value = 50.16666666666667
